$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.094.16'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '1.557.78'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '292.33'
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3987'
$ws.Range("E7").Value = '  +5.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3233'
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.00'
$ws.Range("E9").Value = '  -1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07326'
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.084'
$ws.Range("E11").Value = '  -5.65%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.97'
$ws.Range("E13").Value = '  -7.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.700'
$ws.Range("E14").Value = '  -3.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001140'
$ws.Range("E15").Value = '  +5.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.658'
$ws.Range("E16").Value = '  -1.92%  '
$ws.Range("D17").Value = '1.556.59'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06598'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.87'
$ws.Range("E19").Value = '  -3.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  +0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.314'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.74'
$ws.Range("E22").Value = '  -3.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.30'
$ws.Range("E23").Value = '  -3.88%  '
$ws.Range("D24").Value = '22.106.53'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.362'
$ws.Range("E25").Value = '  +2.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.446'
$ws.Range("E26").Value = '  -6.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '148.75'
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.66'
$ws.Range("E28").Value = '  -3.92%  '
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("D30").Value = '1.730.36'
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.17'
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.015'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.764'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08357'
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.621'
$ws.Range("E35").Value = '  -16.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.118'
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02278'
$ws.Range("E37").Value = '  -4.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06136'
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.151'
$ws.Range("E39").Value = '  -4.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.219'
$ws.Range("E40").Value = '  -2.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2062'
$ws.Range("E41").Value = '  -4.79%  '
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.77'
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5871'
$ws.Range("E44").Value = '  -4.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.766'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.09'
$ws.Range("E46").Value = '  -5.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5616'
$ws.Range("E47").Value = '  -5.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.92'
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.912'
$ws.Range("E49").Value = '  -4.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.141'
$ws.Range("E50").Value = '  -3.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06850'
$ws.Range("E51").Value = '  -3.73%  '
